$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update restated figures on existing rows 173-175 (monthly revision) ---

# Row 173
$ws.Cells.Item(173, 2).Value = 303238    # B173
$ws.Cells.Item(173, 6).Value = 28708     # F173
$ws.Cells.Item(173, 7).Value = 42914     # G173
$ws.Cells.Item(173, 8).Value = 32869     # H173
$ws.Cells.Item(173, 9).Value = 121627    # I173
$ws.Cells.Item(173, 18).Value = 231964   # R173
$ws.Cells.Item(173, 22).Value = 24435    # V173
$ws.Cells.Item(173, 23).Value = 34957    # W173
$ws.Cells.Item(173, 24).Value = 24723    # X173
$ws.Cells.Item(173, 25).Value = 88311    # Y173

# Row 174
$ws.Cells.Item(174, 2).Value = 299480    # B174
$ws.Cells.Item(174, 5).Value = 53456     # E174
$ws.Cells.Item(174, 6).Value = 30953     # F174
$ws.Cells.Item(174, 7).Value = 41406     # G174
$ws.Cells.Item(174, 8).Value = 32671     # H174
$ws.Cells.Item(174, 9).Value = 117672    # I174
$ws.Cells.Item(174, 18).Value = 237004   # R174
$ws.Cells.Item(174, 21).Value = 42953    # U174
$ws.Cells.Item(174, 22).Value = 26804    # V174
$ws.Cells.Item(174, 23).Value = 34298    # W174
$ws.Cells.Item(174, 24).Value = 26050    # X174
$ws.Cells.Item(174, 25).Value = 89255    # Y174

# Row 175
$ws.Cells.Item(175, 2).Value = 286021    # B175
$ws.Cells.Item(175, 4).Value = 19623     # D175
$ws.Cells.Item(175, 5).Value = 41505     # E175
$ws.Cells.Item(175, 6).Value = 29600     # F175
$ws.Cells.Item(175, 7).Value = 42614     # G175
$ws.Cells.Item(175, 8).Value = 32794     # H175
$ws.Cells.Item(175, 9).Value = 117854    # I175
$ws.Cells.Item(175, 18).Value = 226989   # R175
$ws.Cells.Item(175, 20).Value = 14957    # T175
$ws.Cells.Item(175, 21).Value = 32898    # U175
$ws.Cells.Item(175, 22).Value = 25450    # V175
$ws.Cells.Item(175, 23).Value = 35726    # W175
$ws.Cells.Item(175, 24).Value = 26472    # X175
$ws.Cells.Item(175, 25).Value = 89939    # Y175

# --- Append new row 176 for period 01-07-2021 ---

# Column A holds a text label that looks like a date ("01-07-2021"). Assigning
# it directly via .Value triggers Excel's automatic date recognition (it would
# be stored as a date serial number with a date number format), which does not
# match how the rest of column A is stored (plain shared-string text). To get
# a genuine text value we build it as a formula that evaluates to the string,
# then convert that cell to a static value via copy / paste-special-values --
# this mirrors how Excel "types" a value without the live-input date sniffing.
$ws.Cells.Item(176, 1).Formula = '="01-07-2021"'
$ws.Cells.Item(176, 1).Copy()
$ws.Cells.Item(176, 1).PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item(176, 2).Value = 278745    # B176
$ws.Cells.Item(176, 3).Value = 2846      # C176
$ws.Cells.Item(176, 4).Value = 12156     # D176
$ws.Cells.Item(176, 5).Value = 40676     # E176
$ws.Cells.Item(176, 6).Value = 28325     # F176
$ws.Cells.Item(176, 7).Value = 42119     # G176
$ws.Cells.Item(176, 8).Value = 34157     # H176
$ws.Cells.Item(176, 9).Value = 118466    # I176
$ws.Cells.Item(176, 10).Value = 56351    # J176
$ws.Cells.Item(176, 11).Value = 682      # K176
$ws.Cells.Item(176, 12).Value = 3462     # L176
$ws.Cells.Item(176, 13).Value = 7802     # M176
$ws.Cells.Item(176, 14).Value = 3605     # N176
$ws.Cells.Item(176, 15).Value = 6484     # O176
$ws.Cells.Item(176, 16).Value = 6875     # P176
$ws.Cells.Item(176, 17).Value = 27441    # Q176
$ws.Cells.Item(176, 18).Value = 222394   # R176
$ws.Cells.Item(176, 19).Value = 2164     # S176
$ws.Cells.Item(176, 20).Value = 8694     # T176
$ws.Cells.Item(176, 21).Value = 32873    # U176
$ws.Cells.Item(176, 22).Value = 24720    # V176
$ws.Cells.Item(176, 23).Value = 35635    # W176
$ws.Cells.Item(176, 24).Value = 27282    # X176
$ws.Cells.Item(176, 25).Value = 91025    # Y176
